$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header values -------------------------------------------------------
# LANCIO:
$ws.Range("B2").Value = 4569
# PAIA DA PRODURRE:
$ws.Range("B3").Value = 125

# --- Recalculate TOTALE (col F) for the TAGLIO section (rows 7-20) -------
# using the new PAIA DA PRODURRE quantity of 125 (CONS/PA * 125)
$ws.Range("F7").Value = 10.0
$ws.Range("F8").Value = 22.5
$ws.Range("F9").Value = 16.25
$ws.Range("F10").Value = 16.0
$ws.Range("F11").Value = 4.38
$ws.Range("F12").Value = 4.38
$ws.Range("F13").Value = 2.0
$ws.Range("F14").Value = 5.0
$ws.Range("F15").Value = 2.5
$ws.Range("F16").Value = 14.38
$ws.Range("F17").Value = 3.13
$ws.Range("F18").Value = 2.0
$ws.Range("F19").Value = 2.13
$ws.Range("F20").Value = 2.88

# --- Remove the two obsolete ORLATURA rows (LOGO / TIMBRO FOD) -----------
# Deleting row 22 twice removes old rows 22 (LOGO) and 23 (TIMBRO FOD),
# shifting everything below up by two rows.
$ws.Rows(22).Delete()
$ws.Rows(22).Delete()

# --- Recalculate TOTALE (col F) for the ORLATURA section (now rows 22-39)
# using the new PAIA DA PRODURRE quantity of 125 (CONS/PA * 125)
$ws.Range("F22").Value = 1750.0
$ws.Range("F23").Value = 875.0
$ws.Range("F24").Value = 375.0
$ws.Range("F25").Value = 187.5
$ws.Range("F26").Value = 1000.0
$ws.Range("F27").Value = 500.0
$ws.Range("F28").Value = 1.75
$ws.Range("F29").Value = 0.13
$ws.Range("F30").Value = 0.13
$ws.Range("F31").Value = 37.5
$ws.Range("F32").Value = 187.5
$ws.Range("F33").Value = 87.5
$ws.Range("F34").Value = 25.0
$ws.Range("F35").Value = 350.0
$ws.Range("F36").Value = 125.0
$ws.Range("F37").Value = 125.0
$ws.Range("F38").Value = 3000.0
$ws.Range("F39").Value = 1.25

# --- Remove the trailing AUTORIZZAZIONE block (row 44, now row 42) -------
$ws.Rows(42).Delete()
$ws.Rows(42).Delete()
$ws.Rows(42).Delete()

# --- Selection / active cell mirrors the ORLATURA section header ---------
$ws.Range("A21:F21").Select()
